$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.685.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.28%  "

$ws.Range("D3").Value = "'2.442.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.90%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'521.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.92%  "

$ws.Range("D6").Value = "'128.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.82%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").Value = "'0.561"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("D9").Value = "'0.0970"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.16%  "

$ws.Range("D10").Value = "'0.150"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.38%  "

$ws.Range("D11").Value = "'4.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.67%  "

$ws.Range("D12").Value = "'0.320"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.71%  "

$ws.Range("D13").Value = "'2.877.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.96%  "

$ws.Range("D14").Value = "'57.668.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.51%  "

$ws.Range("D15").Value = "'21.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.29%  "

$ws.Range("E16").Value = "  -3.84%  "

$ws.Range("D17").Value = "'2.445.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.91%  "

$ws.Range("D18").Value = "'10.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.99%  "

$ws.Range("D19").Value = "'4.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.91%  "

$ws.Range("D20").Value = "'311.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.97%  "

$ws.Range("D21").Value = "'6.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.64%  "

$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").Value = "'64.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.67%  "

$ws.Range("D24").Value = "'0.400"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.82%  "

$ws.Range("D26").Value = "'2.563.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.79%  "

$ws.Range("E27").Value = "  -3.57%  "

$ws.Range("D28").Value = "'7.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.26%  "

$ws.Range("D29").Value = "'173.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.07%  "

$ws.Range("D30").Value = "'0.0₃0731"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.85%  "

$ws.Range("E31").Value = "  -3.38%  "

$ws.Range("E32").Value = "  -3.65%  "

$ws.Range("E33").Value = "  -9.23%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").Value = "'17.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.88%  "

$ws.Range("E37").Value = "  -7.62%  "

$ws.Range("D38").Value = "'3.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.16%  "

$ws.Range("D39").Value = "'36.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.37%  "

$ws.Range("D40").Value = "'0.795"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.95%  "

$ws.Range("D41").Value = "'1.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.69%  "

$ws.Range("D42").Value = "'3.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.84%  "

$ws.Range("D43").Value = "'0.580"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.85%  "

$ws.Range("D44").Value = "'4.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.03%  "

$ws.Range("D45").Value = "'0.0916"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.88%  "

$ws.Range("D46").Value = "'252.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.02%  "

$ws.Range("D47").Value = "'120.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.44%  "

$ws.Range("D48").Value = "'0.0488"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.73%  "

$ws.Range("E49").Value = "  -3.74%  "

$ws.Range("D50").Value = "'16.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.69%  "

$ws.Range("D51").Value = "'16.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.62%  "
